$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fill in row 22 (day 21) ---
$ws.Cells.Item(22, 1).Value = 21

$ws.Range("B22").NumberFormat = "d-mmm-yy"
$ws.Cells.Item(22, 2).Value = 43851

$ws.Cells.Item(22, 3).Value = "Surah Baqarah, 218 - 229"

$ws.Cells.Item(22, 6).Value = "Women in Islam, Zarrar (ra), Year 336, Battle of Yarmouk"

$content = "h3: Feminism in Islam`np: Some questions that bothers us, made our women go out on strike fighting for their rights.`np.b-left: Why is it not allowed for women to play in open grounds side by side to men?`np.b-left: Why our daughters can not go out in open just like men do?`np.b-left: Why men are given more rights then women?`np.b-left: Why is this life so difficult for women; what is their fault?`np.b-left: Why can" + [char]0x2019 + "t our daughters participate freely in athletics / olympics / sports?`np: Something has impacted our lives. We were not always living our lives this way. It happened very recently that women felt insecure and went out fighting to win security. Men became bolder and angry to beat women to their knees, as and when they found a women raising her voice. Happening in areas of low literacy, in Asia, 3rd world countries, backward areas. `np: USA, Europe, modern world still allows women to live their freely. A man raising his voice, hand on a women is made to undergo a complete cycle of investigation. Not happening around here. Man raising his voice in this part of the world is Islamic-ally accepted and allowed. `np: People who follow West flourish. People who do not, decline. Does West follow Islam? Yes to some extent they follow, but Islam does not allow our girls to wear short dresses and play football in public.. This confusion leads our daughters astray from our religion. They go all out, fighting for their rights and some times they become Qandeel Baloch and at other times just a snubbed voice in the corner of a home.`np: Girls are given equal rights in Islam:-`nquote: And due to the wives is similar to what is expected of them, according to what is reasonable. But the men have a degree over them [in responsibility and authority]. <br> - Surah Baqarah verse 228`np: Girls are given rights but they are not equal in responsibility. Something to understand here, the lynchpin, if they are equal in rights then why do Men have a degree over them?`np: Men are physically stronger, responsible to take burden of our homes, stable in decision making. The in built functionality of men can not be challenged. Women are equal to them in rights. But the functionality of women is different. They are impulsive and emotional. `np: <b>Now, is this the case with all women? </b>`np: No. All women are not emotional and impulsive. Just like all men are not responsible. All cars are not reliable. Some break the first day they come on the road. Cars that are designed and built in-houses are weaker compared to the ones built by Toyota. Men that come on the road of responsibility out of pampered houses break the first day. `np: Inherently, genetically, a man feels better if he leads. A women feels better if she is loved. They are opposite sides of the magnet attracting each other. `np: Allah (swt) is all wise and all knower. He plans things His way. My justification can be wrong but the incident below is not wrong.`nh3: Khawlah bint al-Azwar`np: Khawlah bint al-Azwar was a Muslim warrior during the life of the Islamic prophet Muhammad and later a military leader. She has been described as one of the greatest female military leaders in history and was once compared with Khalid Bin Walid by the opponents in the battlefield. She was the sister of Dhiraar bin Al-Azwar, the soldier and commander of the Rashidun army during the 7th century Muslim conquest. Born sometime in the seventh century, daughter of Malik or Taken bin Awake, one of the chiefs of the Banu Assad tribe, Khawlah was well known for her leadership in battles of the Muslim conquests in parts of what are today Syria, Jordan, and Palestine. She fought side by side with her brother Zarrar (ra) in many battles, including the decisive Battle of Yarmouk in 636 against the Byzantine empire. On the 4th day of the battle she led a group of women against the Byzantine army and defeated its chief commander and later was wounded during her fight with a Greek soldier.`np.note: <a href=" + [char]0x201C + "https://en.wikipedia.org/wiki/Khawlah_bint_al-Azwar" + [char]0x201D + ">Click to go to her wikipedia page</a>`np: Girls can fight the battles, lead them riding horses wearing armour. They can play football with all their energies, get injured, fall and rise. They can go do athletics. It is just fine for them to roam in public <b>covered</b>. They have same books to teach themselves from that exist for others to learn from. They can excel in every corner of this life that they want to. `nh3: Let" + [char]0x2019 + "s stay put..`np: Staying put to Islam does not limit us, it protects us from cringing complexities of society we live in. `np: Islam is a wonderful religion. It frees every bit of us. No dowry, no jewlery, no prison cells. You divorce you don" + [char]0x2019 + "t want to live together. You do not want to go to Jirga / elderly in your tribe to get divorce approved, do not go. `np: You just get a divorce if it does not work. It is not as worse as " + [char]0x201C + "telling a lie" + [char]0x201D + " in public. Lies have become sober, normal and practical solution to our lives. `np: And divorce? It has become the <b>" + [char]0x201D + "least liked ruling to Allah (swt)" + [char]0x201D + "</b>?`np: What Allah (swt) hates we are holding tight to it. And what Allah (swt) allows, we use it to break our daughters.`nquote: And if he has divorced her [for the third time], then she is not lawful to him afterward until [after] she marries a husband other than him. And if the latter husband divorces her [or dies], there is no blame upon the woman and her former husband for returning to each other if they think that they can keep [within] the limits of Allah . These are the limits of Allah, which He makes clear to a people who know. <br> - Surah Baqarah verse 230`np.note: Guys, if you like this project. Please follow this project's page on twitter. <a href=`"https://twitter.com/zakatlists`">Click here to go to the twitter page</a>."

$ws.Cells.Item(22, 4).Value = $content

$ws.Cells.Item(22, 5).Value = "Qasim Ali"

$ws.Rows.Item(22).RowHeight = 409.6

# --- Update the sheet view / selection to the newly active row ---
$ws.Activate()
[void]$ws.Range("D22").Select()
